$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values that changed (rows 2-16, columns B-F)
$ws.Range("B2").Value = "NSE:ADANIENT"
$ws.Range("C2").Value = "NSE:AAKASH"
$ws.Range("D2").Value = "NSE:APLAPOLLO"
$ws.Range("E2").Value = "NSE:GODREJPROP"
$ws.Range("B3").Value = "NSE:ASTERDM"
$ws.Range("C3").Value = "NSE:AMBIKCO"
$ws.Range("D3").Value = "NSE:KEI"
$ws.Range("E3").Value = "NSE:NHPC"
$ws.Range("F3").Value = "NSE:CHOLAFIN"
$ws.Range("B4").Value = "NSE:BAJAJHLDNG"
$ws.Range("C4").Value = "NSE:APEX"
$ws.Range("E4").Value = "NSE:OIL"
$ws.Range("F4").Value = "NSE:CIPLA"
$ws.Range("B5").Value = "NSE:CAMPUS"
$ws.Range("C5").Value = "NSE:BFINVEST"
$ws.Range("E5").Value = "NSE:PIDILITIND"
$ws.Range("B6").Value = "NSE:CAMS"
$ws.Range("C6").Value = "NSE:BIRLAMONEY"
$ws.Range("B7").Value = "NSE:HEALTHY"
$ws.Range("C7").Value = "NSE:BLKASHYAP"
$ws.Range("B8").Value = "NSE:KPIGREEN"
$ws.Range("C8").Value = "NSE:BPL"
$ws.Range("B9").Value = "NSE:LGHL"
$ws.Range("C9").Value = "NSE:DCMSRIND"
$ws.Range("B10").Value = "NSE:NGIL"
$ws.Range("C10").Value = "NSE:GLOBUSSPR"
$ws.Range("B11").Value = "NSE:PHARMABEES"
$ws.Range("C11").Value = "NSE:KARMAENG"
$ws.Range("B12").Value = "NSE:RADICO"
$ws.Range("C12").Value = "NSE:KOHINOOR"
$ws.Range("C13").Value = "NSE:MARINE"
$ws.Range("C14").Value = "NSE:OMINFRAL"
$ws.Range("C15").Value = "NSE:ONWARDTEC"
$ws.Range("C16").Value = "NSE:SADHNANIQ"

# Clear cells that became empty
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("B13").ClearContents()

# Remove now-unused rows 17-22 (table shrank from 20 tickers to 15)
$ws.Range("A17:A22").EntireRow.Delete() | Out-Null

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
